$d = $word.ActiveDocument

# Add the new "NormalLine" character style (wdStyleTypeCharacter = 2)
$style = $d.Styles.Add("NormalLine", 2)
$style.Priority = 1
$style.QuickStyle = $true

$style.Font.Name = "Calibri"
$style.Font.Size = 11
$style.Font.SizeBi = 12
$style.Font.Color = -16777216
